$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row = 2; D = 44181; J = 38; K = 26000; L = 26000; M = 26000; N = "$/malla 25 kilos"; O = "Región Metropolitana"; P = 1040 }
  @{ Row = 3; D = 44412; J = 35; K = 24000; L = 24000; M = 24000; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 960 }
  @{ Row = 4; D = 44159; J = 35; K = 22000; L = 22000; M = 22000; N = "$/malla 25 kilos"; O = "Provincia de Quillota"; P = 880 }
  @{ Row = 5; D = 44406; J = 35; K = 32000; L = 32000; M = 32000; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 1280 }
  @{ Row = 6; D = 44253; J = 38; K = 18000; L = 18000; M = 18000; N = "$/saco 25 kilos"; O = "Provincia de Talca"; P = 720 }
  @{ Row = 7; D = 44250; J = 38; K = 18000; L = 18000; M = 18000; N = "$/malla 25 kilos"; O = "Provincia de Talca"; P = 720 }
  @{ Row = 8; D = 44399; J = 38; K = 33000; L = 33000; M = 33000; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 1320 }
  @{ Row = 9; D = 44452; J = 70; K = 31000; L = 32000; M = 31500; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 1260 }
  @{ Row = 10; D = 44161; J = 35; K = 21000; L = 21000; M = 21000; N = "$/saco 25 kilos"; O = "Provincia de Quillota"; P = 840 }
  @{ Row = 11; D = 44453; J = 73; K = 21000; L = 22000; M = 21521; N = "$/saco 25 kilos"; O = "Provincia de Limarí"; P = 861 }
  @{ Row = 12; D = 44448; J = 45; K = 32000; L = 32000; M = 32000; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 1280 }
  @{ Row = 13; D = 44370; J = 45; K = 32000; L = 32000; M = 32000; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 1280 }
  @{ Row = 14; D = 44165; J = 45; K = 22000; L = 22000; M = 22000; N = "$/saco 25 kilos"; O = "Provincia de Quillota"; P = 880 }
  @{ Row = 15; D = 44343; J = 40; K = 28000; L = 28000; M = 28000; N = "$/saco 25 kilos"; O = "Provincia de Limarí"; P = 1120 }
  @{ Row = 16; D = 44252; J = 40; K = 18000; L = 19000; M = 18625; N = "$/malla 25 kilos"; O = "Provincia de Talca"; P = 745 }
  @{ Row = 17; D = 44372; J = 50; K = 33000; L = 34000; M = 33500; N = "$/saco 25 kilos"; O = "Provincia de Limarí"; P = 1340 }
  @{ Row = 18; D = 44160; J = 35; K = 21000; L = 21000; M = 21000; N = "$/saco 25 kilos"; O = "Provincia de Quillota"; P = 840 }
  @{ Row = 19; D = 44365; J = 70; K = 22000; L = 23000; M = 22500; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 900 }
  @{ Row = 20; D = 44162; J = 35; K = 17000; L = 17000; M = 17000; N = "$/saco 25 kilos"; O = "Provincia de Quillota"; P = 680 }
  @{ Row = 21; D = 44410; J = 35; K = 34000; L = 34000; M = 34000; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 1360 }
  @{ Row = 22; D = 44411; J = 35; K = 34000; L = 34000; M = 34000; N = "$/malla 25 kilos"; O = "Provincia de Limarí"; P = 1360 }
  @{ Row = 23; D = 44376; J = 38; K = 27000; L = 27000; M = 27000; N = "$/saco 25 kilos"; O = "Provincia de Limarí"; P = 1080 }
)

foreach ($row in $rows) {
  $r = $row.Row
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 10).Value = $row.J
  $ws.Cells.Item($r, 11).Value = $row.K
  $ws.Cells.Item($r, 12).Value = $row.L
  $ws.Cells.Item($r, 13).Value = $row.M
  $ws.Cells.Item($r, 14).Value = $row.N
  $ws.Cells.Item($r, 15).Value = $row.O
  $ws.Cells.Item($r, 16).Value = $row.P
}
